$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the two hyperlinks that were on the old P2 (thumbnail_img) and
#    Q2 (photos) cells - the new layout has no hyperlinks at all.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Relocate the trailing "meta_title / meta_description / thumbnail_img /
#    photos" block from columns N:Q out to the new end of the table
#    (columns W:Z), since columns N:Q are being reused for new fields.
#    Copy header text + formatting first (before anything is overwritten).
# ---------------------------------------------------------------------------
$ws.Range("N1:Q1").Copy() | Out-Null
$ws.Range("W1:Z1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("W1").Value = "meta_title"
$ws.Range("X1").Value = "meta_description"
$ws.Range("Y1").Value = "thumbnail_img"
$ws.Range("Z1").Value = "photos"

# ---------------------------------------------------------------------------
# 3. Overwrite the (now freed) header cells N1:V1 with the new field names,
#    carrying over the bold header formatting used by the rest of row 1.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("N1:V1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("N1").Value = "returnable"
$ws.Range("O1").Value = "returnable_days"
$ws.Range("P1").Value = "discount"
$ws.Range("Q1").Value = "discount_type"
$ws.Range("R1").Value = "earn_point"
$ws.Range("S1").Value = "expiry_month"
$ws.Range("T1").Value = "expiry_year"
$ws.Range("U1").Value = "batch_number"
$ws.Range("V1").Value = "hsn_code"

# ---------------------------------------------------------------------------
# 4. Update the demo data row (row 2) for the new "Demo Product 23" sample.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = "Demo Product 23"
$ws.Range("B2").Value = "Demo Product Description 23"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 1

# video_provider / video_link no longer populated for this sample
$ws.Range("E2").ClearContents()
$ws.Range("F2").ClearContents()

$ws.Range("G2").Value = "Demo, Baby"
$ws.Range("H2").Value = 200
$ws.Range("I2").Value = 120
$ws.Range("J2").Value = "pc"
$ws.Range("K2").Value = "demo-product-23"
$ws.Range("L2").Value = 35
$ws.Range("M2").Value = "SKU-23"

$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 12
$ws.Range("P2").Value = 10
$ws.Range("Q2").Value = "percent"
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = "Feb"
$ws.Range("T2").Value = 2023
$ws.Range("U2").Value = "BT-OCT22"
$ws.Range("V2").Value = "HSN3029"

# old thumbnail/photo hyperlink text is gone - P2/Q2 lose the Hyperlink look,
# and the relocated meta/thumbnail/photo columns (W2:Z2) are left blank but
# Y2/Z2 (thumbnail_img/photos) keep the Hyperlink-esque styling of the old
# P2/Q2 cells even though empty.
$ws.Range("P2").Style = "Normal"
$ws.Range("Q2").Style = "Normal"
$ws.Range("W2").ClearContents()
$ws.Range("X2").ClearContents()
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("Y2").Style = "Hyperlink"
$ws.Range("Z2").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5. Widen the newly-inserted columns so the new headers/content are not
#    squeezed into the sheet's default column width.
# ---------------------------------------------------------------------------
$ws.Range("N1:V1").ColumnWidth = 15.9
$ws.Range("W1").ColumnWidth = 11.9
$ws.Range("X1").ColumnWidth = 18.3
$ws.Range("Y1").ColumnWidth = 126.4

# ---------------------------------------------------------------------------
# 6. Selection cosmetics to match the saved workbook (active cell moved to
#    the new last column of row 2).
# ---------------------------------------------------------------------------
$ws.Range("AB2").Select() | Out-Null
